$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 5609.6313
$ws.Range("I28").Value = 6230.5293
$ws.Range("J28").Value = 332
$ws.Range("K28").Value = 6230.5293
$ws.Range("L28").Value = 332
$ws.Range("M28").Value = -5745.5293
$ws.Range("N28").Value = -1302

$ws.Range("H52").Value = 5000
$ws.Range("I52").Value = 9000
$ws.Range("K52").Value = 27000
$ws.Range("M52").Value = -26840

$ws.Range("H74").Value = 25800.6
$ws.Range("I74").Value = 25800.6
$ws.Range("K74").Value = 25800.6
$ws.Range("M74").Value = -24864.6

$ws.Range("H77").Value = 25800.6
$ws.Range("I77").Value = 25800.6
$ws.Range("K77").Value = 129003
$ws.Range("M77").Value = -124323

$ws.Range("H86").Value = 2301.15
$ws.Range("J86").Value = 1927.625
$ws.Range("L86").Value = 1927.625
$ws.Range("N86").Value = -4173.625

$ws.Range("H89").Value = 2301.15
$ws.Range("J89").Value = 1927.625
$ws.Range("L89").Value = 9638.125
$ws.Range("N89").Value = -20870.125

$ws.Range("H112").Value = 1968.7273
$ws.Range("I112").Value = 1292.7142
$ws.Range("J112").Value = 3151.75
$ws.Range("K112").Value = 3878.1426
$ws.Range("L112").Value = 9455.25
$ws.Range("M112").Value = -2770.1426
$ws.Range("N112").Value = -11671.25

$ws.Range("H113").Value = 3957.3333
$ws.Range("I113").Value = 3957.3333
$ws.Range("K113").Value = 3957.3333
$ws.Range("M113").Value = -703.3332999999998

$ws.Range("H125").Value = 67666.336
$ws.Range("I125").Value = 1499
$ws.Range("K125").Value = 13491
$ws.Range("M125").Value = -11031

$ws.Range("H137").Value = 1307.1708
$ws.Range("J137").Value = 780
$ws.Range("L137").Value = 2340
$ws.Range("N137").Value = -7440

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2703.8667
$ws.Range("I2").Value = 1996
$ws.Range("K2").Value = 1996
$ws.Range("M2").Value = -1883

$ws.Range("H32").Value = 4461.745
$ws.Range("I32").Value = 3494.4792
$ws.Range("K32").Value = 3494.4792
$ws.Range("M32").Value = -3207.4792

$ws.Range("H33").Value = 73351.664
$ws.Range("I33").Value = 60013
$ws.Range("J33").Value = 100029
$ws.Range("K33").Value = 60013
$ws.Range("L33").Value = 100029
$ws.Range("M33").Value = -59684
$ws.Range("N33").Value = -100687

$ws.Range("H102").Value = 3026.8157
$ws.Range("I102").Value = 2016.6552
$ws.Range("J102").Value = 6281.778
$ws.Range("K102").Value = 2016.6552
$ws.Range("L102").Value = 6281.778
$ws.Range("M102").Value = -394.6551999999999
$ws.Range("N102").Value = -9525.778

$ws.Range("H116").Value = 2703.8667
$ws.Range("I116").Value = 1996
$ws.Range("K116").Value = 1996
$ws.Range("M116").Value = 298

$ws.Range("H122").Value = 2714.7273
$ws.Range("I122").Value = 2692
$ws.Range("K122").Value = 8076
$ws.Range("M122").Value = -5626

$ws.Range("H132").Value = 1569.48
$ws.Range("I132").Value = 1370.8096
$ws.Range("K132").Value = 4112.4288
$ws.Range("M132").Value = -1582.4288

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2703.8667
$ws.Range("I3").Value = 1996
$ws.Range("K3").Value = 1996
$ws.Range("M3").Value = -1882

$ws.Range("H28").Value = 30000
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").Value = ""

$ws.Range("H40").Value = 39999.8
$ws.Range("J40").Value = 39999.8
$ws.Range("L40").Value = 39999.8
$ws.Range("N40").Value = -40529.8

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = ""

$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").Value = ""

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = ""

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 968
$ws.Range("I22").Value = 898.25
$ws.Range("K22").Value = 898.25
$ws.Range("M22").Value = -548.25

$ws.Range("H58").Value = 1265.1
$ws.Range("I58").Value = 1407.2858
$ws.Range("J58").Value = 933.3333
$ws.Range("K58").Value = 1407.2858
$ws.Range("L58").Value = 933.3333
$ws.Range("M58").Value = -1204.2858
$ws.Range("N58").Value = -1339.3333

$ws.Range("H122").Value = 924.875
$ws.Range("I122").Value = 966.5
$ws.Range("J122").Value = 800
$ws.Range("K122").Value = 2899.5
$ws.Range("L122").Value = 2400
$ws.Range("M122").Value = -449.5
$ws.Range("N122").Value = -7300

$ws.Range("H136").Value = 1265.1
$ws.Range("I136").Value = 1407.2858
$ws.Range("J136").Value = 933.3333
$ws.Range("K136").Value = 4221.857400000001
$ws.Range("L136").Value = 2799.9999
$ws.Range("M136").Value = -1671.857400000001
$ws.Range("N136").Value = -7899.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").Value = ""

$ws.Range("H44").Value = 1468.1666
$ws.Range("I44").Value = 1638.2858
$ws.Range("J44").Value = 1230
$ws.Range("K44").Value = 4914.857400000001
$ws.Range("L44").Value = 3690
$ws.Range("M44").Value = -4516.857400000001
$ws.Range("N44").Value = -4486

$ws.Range("H51").Value = 1049.6666
$ws.Range("I51").Value = 599.3333
$ws.Range("K51").Value = 1797.9999
$ws.Range("M51").Value = -1337.9999

$ws.Range("H86").Value = 966.75
$ws.Range("J86").Value = 678
$ws.Range("L86").Value = 2034
$ws.Range("N86").Value = -4406

$ws.Range("H89").Value = 966.75
$ws.Range("J89").Value = 678
$ws.Range("L89").Value = 6102
$ws.Range("N89").Value = -17958

$ws.Range("H120").Value = 37772.363
$ws.Range("I120").Value = 16799.2
$ws.Range("J120").Value = 55250
$ws.Range("K120").Value = 50397.60000000001
$ws.Range("L120").Value = 165750
$ws.Range("M120").Value = -45559.60000000001
$ws.Range("N120").Value = -175426

$ws.Range("H129").Value = 2874
$ws.Range("J129").Value = 4273.6665
$ws.Range("L129").Value = 12820.9995
$ws.Range("N129").Value = -22820.9995

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 1676833.4
$ws.Range("J20").Value = 12200
$ws.Range("L20").Value = 12200
$ws.Range("N20").Value = -12690

$ws.Range("H45").Value = 69420
$ws.Range("J45").Value = 69420
$ws.Range("L45").Value = 69420
$ws.Range("N45").Value = -70538

$ws.Range("H51").Value = 57844
$ws.Range("J51").Value = 75320.8
$ws.Range("L51").Value = 75320.8
$ws.Range("N51").Value = -76338.8

$ws.Range("H102").Value = 8155.375
$ws.Range("I102").Value = 4798.8
$ws.Range("J102").Value = 13749.667
$ws.Range("K102").Value = 4798.8
$ws.Range("L102").Value = 13749.667
$ws.Range("M102").Value = -3176.8
$ws.Range("N102").Value = -16993.667

$ws.Range("H107").Value = 1932.2858
$ws.Range("I107").Value = 152
$ws.Range("J107").Value = 2069.2307
$ws.Range("K107").Value = 152
$ws.Range("L107").Value = 2069.2307
$ws.Range("M107").Value = 1768
$ws.Range("N107").Value = -5909.2307

$ws.Range("H132").Value = 2873.3333
$ws.Range("I132").Value = 2873.3333
$ws.Range("K132").Value = 8619.999899999999
$ws.Range("M132").Value = -6089.999899999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6833.3335
$ws.Range("I7").Value = 5850
$ws.Range("K7").Value = 5850
$ws.Range("M7").Value = -5738

$ws.Range("H40").Value = 4350.3447
$ws.Range("I40").Value = 2874.2856
$ws.Range("J40").Value = 8225
$ws.Range("K40").Value = 2874.2856
$ws.Range("L40").Value = 8225
$ws.Range("M40").Value = -2738.2856
$ws.Range("N40").Value = -8497

$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").Value = ""

$ws.Range("H46").Value = 1541.0408
$ws.Range("I46").Value = 1088.2
$ws.Range("J46").Value = 1853.3448
$ws.Range("K46").Value = 1088.2
$ws.Range("L46").Value = 1853.3448
$ws.Range("M46").Value = -900.2
$ws.Range("N46").Value = -2229.3448

$ws.Range("H122").Value = 3724.2856
$ws.Range("I122").Value = 4302.8335
$ws.Range("J122").Value = 2952.889
$ws.Range("K122").Value = 12908.5005
$ws.Range("L122").Value = 8858.667000000001
$ws.Range("M122").Value = -10458.5005
$ws.Range("N122").Value = -13758.667

$ws.Range("H126").Value = 6833.3335
$ws.Range("I126").Value = 5850
$ws.Range("K126").Value = 17550
$ws.Range("M126").Value = -15080

$ws.Range("H132").Value = 3834.4075
$ws.Range("I132").Value = 3660.2273
$ws.Range("K132").Value = 10980.6819
$ws.Range("M132").Value = -8450.6819
